# PO Clean up - updated the existing TC
# - Insert a new "GLAccount" worksheet right before "JournalEntries"
# - Populate it with an "Account Name" header and a sample GL account value
# - Update the selected cell / active sheet state on a couple of other sheets

$wb = $excel.ActiveWorkbook

# ObjectName sheet: selection moves from F1:G1048576 to H6 (no longer the active tab)
$objectName = $wb.Worksheets.Item("ObjectName")
$objectName.Range("H6").Select()

# JournalEntriesD sheet: selection moves from H15 to A2
$journalEntriesD = $wb.Worksheets.Item("JournalEntriesD")
$journalEntriesD.Range("A2").Select()

# Add the new GLAccount sheet immediately before JournalEntries
$journalEntries = $wb.Worksheets.Item("JournalEntries")
$glAccount = $wb.Worksheets.Add($journalEntries)
$glAccount.Name = "GLAccount"

$glAccount.Range("A1").Value = "Account Name"
$glAccount.Range("A2").Value = "2000 (Accounts Payable)"
$glAccount.Columns.Item(1).ColumnWidth = 22.7109375

# Leave GLAccount as the active sheet with A2 selected
$glAccount.Range("A2").Select()
